$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows 180-186: extracted from 10.1016/j.optlastec.2019.01.009 ---

# Column B (alloy composition) - write order matters for shared-string table order.
$ws.Range("B180").Value = "(AlMoCrFe)24.925 V0.3"
$ws.Range("B181").Value = "(AlMoCrFe)23.125 V7.5"
$ws.Range("B182").Value = "(AlMoCrFe)22.375 V10.5"
$ws.Range("B183").Value = "(AlMoCrFe)21.05 V15.8"
$ws.Range("B184").Value = "AlMoCrFeV"
$ws.Range("B185").Value = "(AlMoCrFe)24.925 V0.3"
$ws.Range("B186").Value = "AlMoCrFeV"

# Column C (phase) - re-uses existing shared string "BCC"
$ws.Range("C180:C186").Value = "BCC"

# Columns D & E for the as-built (FGM) rows 180-184
$ws.Range("D180:D184").Value = "LENS(AM)"
$ws.Range("E180:E184").Value = "Laser Engineered Net Shaping (LENS) AM; 1500 W Nd: YAG laser with oxygen content under 10ppm; FGM"

# Columns D & E for the annealed rows 185-186
$ws.Range("D185:D186").Value = "LENS(AM)+A"
$ws.Range("E185:E186").Value = "Laser Engineered Net Shaping (LENS) AM; 1500 W Nd: YAG laser with oxygen content under 10ppm; FGM; annealed at 1373K for 30min in Ar"

# Column F (property measured) - re-uses existing shared string "hardness"
$ws.Range("F180:F186").Value = "hardness"

# Column G (data type) - re-uses existing shared string "EXP"
$ws.Range("G180:G186").Value = "EXP"

# Column I (temperature, K)
$ws.Range("I180:I186").Value = 298

# Column L (units) - re-uses existing shared string "Pa"
$ws.Range("L180:L184").Value = "Pa"

# Column M (temperature label) - re-uses existing shared string "T2"
$ws.Range("M180:M184").Value = "T2"

# Column N (DOI reference) - first write at N180 introduces the new shared string.
$ws.Range("N180:N186").Value = "10.1016/j.optlastec.2019.01.009"

# Column P / Q (raw hardness / uncertainty values feeding the formulas below)
$ws.Range("P180").Value = 485
$ws.Range("Q180").Value = 12
$ws.Range("P181").Value = 521
$ws.Range("Q181").Value = 7
$ws.Range("P182").Value = 542
$ws.Range("Q182").Value = 18
$ws.Range("P183").Value = 558
$ws.Range("Q183").Value = 14
$ws.Range("P184").Value = 581
$ws.Range("Q184").Value = 21
$ws.Range("P185").Value = 435
$ws.Range("P186").Value = 536

# Columns J & K: shared formula converting hardness to Pa
$ws.Range("J180:K186").Formula = "=P180*9807000"
# Rows 185-186 only report the converted hardness (J); no uncertainty (K) was given.
$ws.Range("K185:K186").ClearContents()

# --- View state: leave the selection where the editor finished working ---
$win = $excel.ActiveWindow
$win.ScrollRow = 162
[void]$ws.Range("N191").Select()
